$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.607.28"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "1.850.20"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.75"
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5251"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3250"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06814"
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.96"
$ws.Range("E10").Value = "  +0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7831"
$ws.Range("E11").Value = "  +1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07797"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").Value = "1.858.51"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.62"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.031"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.02"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007991"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").Value = "26.618.20"
$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.643"
$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.496"
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.030"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.92"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.191"
$ws.Range("E25").Value = "  -7.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.682"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.08"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.13"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.206"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.128"
$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08736"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04841"
$ws.Range("E32").Value = "  +0.88%  "

$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7235"
$ws.Range("E34").Value = "  +5.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.884"
$ws.Range("E35").Value = "  +1.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.113"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.281"
$ws.Range("E37").Value = "  +3.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01799"
$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4879"
$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9035"
$ws.Range("E40").Value = "  +0.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.02"
$ws.Range("E41").Value = "  -1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.001"
$ws.Range("E42").Value = "  -3.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.682"
$ws.Range("E44").Value = "  -1.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4219"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05886"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.028"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1238"
$ws.Range("E48").Value = "  -2.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.10"
$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8896"
$ws.Range("E50").Value = "  +3.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.11"
$ws.Range("E51").Value = "  +1.07%  "
